$d = $word.ActiveDocument
$d.Content.Find.Execute("repredentantes", $true, $false, $false, $false, $false, $true, 1, $false, "representantes", 2)
